# Auto-generated edit script
# Updates Leve profit/price calculation columns (H-N) across all 8 sheets
# per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 187.23077
$ws.Range("I5").Value = 184.33333
$ws.Range("J5").Value = 193.75
$ws.Range("K5").Value = 184.33333
$ws.Range("L5").Value = 193.75
$ws.Range("M5").Value = -69.33332999999999
$ws.Range("N5").Value = -423.75
$ws.Range("H40").Value = 4246.1304
$ws.Range("I40").Value = 4382.8
$ws.Range("J40").Value = 4141
$ws.Range("K40").Value = 4382.8
$ws.Range("L40").Value = 4141
$ws.Range("M40").Value = -4207.8
$ws.Range("N40").Value = -4491
$ws.Range("H70").Value = 2999.1538
$ws.Range("I70").Value = 2998.9
$ws.Range("K70").Value = 8996.700000000001
$ws.Range("M70").Value = -8726.700000000001
$ws.Range("H73").Value = 2999.1538
$ws.Range("I73").Value = 2998.9
$ws.Range("K73").Value = 8996.700000000001
$ws.Range("M73").Value = -8060.700000000001
$ws.Range("H82").Value = 500
$ws.Range("I82").Value = 500
$ws.Range("K82").Value = 1500
$ws.Range("M82").Value = -1094
$ws.Range("H85").Value = 500
$ws.Range("I85").Value = 500
$ws.Range("K85").Value = 1500
$ws.Range("M85").Value = -96
$ws.Range("H97").Value = 1106.6666
$ws.Range("J97").Value = 1106.6666
$ws.Range("L97").Value = 3319.9998
$ws.Range("N97").Value = -4311.9998
$ws.Range("H98").Value = 603.5714
$ws.Range("I98").Value = 523.0909
$ws.Range("K98").Value = 523.0909
$ws.Range("M98").Value = 974.9091
$ws.Range("H99").Value = 1495.8572
$ws.Range("J99").Value = 3193.3333
$ws.Range("L99").Value = 9579.999899999999
$ws.Range("N99").Value = -12575.9999
$ws.Range("H101").Value = 225
$ws.Range("I101").Value = 225
$ws.Range("K101").Value = 675
$ws.Range("M101").Value = 947
$ws.Range("H107").Value = 418.72726
$ws.Range("I107").Value = 254.6
$ws.Range("J107").Value = 555.5
$ws.Range("K107").Value = 254.6
$ws.Range("L107").Value = 555.5
$ws.Range("M107").Value = 1665.4
$ws.Range("N107").Value = -4395.5
$ws.Range("H122").Value = 603.5714
$ws.Range("I122").Value = 523.0909
$ws.Range("K122").Value = 1569.2727
$ws.Range("M122").Value = 880.7273
$ws.Range("H135").Value = 832.0833
$ws.Range("I135").Value = 766.35297
$ws.Range("K135").Value = 6897.17673
$ws.Range("M135").Value = -4362.17673
$ws.Range("H138").Value = 3099.0833
$ws.Range("I138").Value = 1084.8572
$ws.Range("J138").Value = 5919
$ws.Range("K138").Value = 3254.5716
$ws.Range("L138").Value = 17757
$ws.Range("M138").Value = 1885.4284
$ws.Range("N138").Value = -28037
$ws.Range("H141").Value = 1175.9
$ws.Range("I141").Value = 1175.9
$ws.Range("K141").Value = 3527.7
$ws.Range("M141").Value = 1652.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3813.0386
$ws.Range("I32").Value = 733.5909
$ws.Range("K32").Value = 733.5909
$ws.Range("M32").Value = -446.5909
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 163
$ws.Range("I80").Value = 158
$ws.Range("J80").Value = 173
$ws.Range("K80").Value = 158
$ws.Range("L80").Value = 173
$ws.Range("M80").Value = 840
$ws.Range("N80").Value = -2169
$ws.Range("H83").Value = 163
$ws.Range("I83").Value = 158
$ws.Range("J83").Value = 173
$ws.Range("K83").Value = 790
$ws.Range("L83").Value = 865
$ws.Range("M83").Value = 4202
$ws.Range("N83").Value = -10849
$ws.Range("H112").Value = 47749.75
$ws.Range("J112").Value = 47749.75
$ws.Range("L112").Value = 47749.75
$ws.Range("N112").Value = -50703.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 36841.5
$ws.Range("J18").Value = 36841.5
$ws.Range("L18").Value = 36841.5
$ws.Range("N18").Value = -37301.5
$ws.Range("H31").Value = 5215.78
$ws.Range("I31").Value = 3574.9375
$ws.Range("K31").Value = 3574.9375
$ws.Range("M31").Value = -3279.9375
$ws.Range("H34").Value = 5215.78
$ws.Range("I34").Value = 3574.9375
$ws.Range("K34").Value = 3574.9375
$ws.Range("M34").Value = -3372.9375
$ws.Range("H107").Value = 748
$ws.Range("J107").Value = 999
$ws.Range("L107").Value = 999
$ws.Range("N107").Value = -4839

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 649.8
$ws.Range("I25").Value = 225
$ws.Range("J25").Value = 933
$ws.Range("K25").Value = 675
$ws.Range("L25").Value = 2799
$ws.Range("M25").Value = -506
$ws.Range("N25").Value = -3137
$ws.Range("H30").Value = 649.8
$ws.Range("I30").Value = 225
$ws.Range("J30").Value = 933
$ws.Range("K30").Value = 675
$ws.Range("L30").Value = 2799
$ws.Range("M30").Value = -573
$ws.Range("N30").Value = -3003
$ws.Range("H109").Value = 707.8889
$ws.Range("I109").Value = 707.8889
$ws.Range("K109").Value = 2123.6667
$ws.Range("M109").Value = -1083.6667
$ws.Range("H132").Value = 3000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1849.4736
$ws.Range("I102").Value = 1841.1111
$ws.Range("K102").Value = 1841.1111
$ws.Range("M102").Value = -219.1111000000001
$ws.Range("H107").Value = 206.66667
$ws.Range("I107").Value = 206.66667
$ws.Range("K107").Value = 206.66667
$ws.Range("M107").Value = 1713.33333
$ws.Range("H113").Value = 2800
$ws.Range("J113").Value = 1199.5
$ws.Range("L113").Value = 1199.5
$ws.Range("N113").Value = -5539.5
$ws.Range("H126").Value = 3099.8
$ws.Range("I126").Value = 3099.8
$ws.Range("K126").Value = 9299.400000000001
$ws.Range("M126").Value = -6829.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 15000
$ws.Range("I5").Value = 15000
$ws.Range("J5").Value = 15000
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = -14887
$ws.Range("N5").Value = -15226
$ws.Range("H22").Value = 1768.1818
$ws.Range("J22").Value = 2062.5
$ws.Range("L22").Value = 2062.5
$ws.Range("N22").Value = -2652.5
$ws.Range("H24").Value = 10007
$ws.Range("J24").Value = 10007
$ws.Range("L24").Value = 10007
$ws.Range("N24").Value = -10693
$ws.Range("H27").Value = 1768.1818
$ws.Range("J27").Value = 2062.5
$ws.Range("L27").Value = 2062.5
$ws.Range("N27").Value = -2276.5
$ws.Range("H46").Value = 1924.6
$ws.Range("I46").Value = 6228
$ws.Range("J46").Value = 848.75
$ws.Range("K46").Value = 6228
$ws.Range("L46").Value = 848.75
$ws.Range("M46").Value = -6040
$ws.Range("N46").Value = -1224.75
$ws.Range("H55").Value = 599.2273
$ws.Range("I55").Value = 655.9474
$ws.Range("J55").Value = 240
$ws.Range("K55").Value = 655.9474
$ws.Range("L55").Value = 240
$ws.Range("M55").Value = -482.9474
$ws.Range("N55").Value = -586

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 13837.333
$ws.Range("J20").Value = 13837.333
$ws.Range("L20").Value = 13837.333
$ws.Range("N20").Value = -14317.333
$ws.Range("H21").Value = 2525000
$ws.Range("I21").Value = 2525000
$ws.Range("K21").Value = 2525000
$ws.Range("M21").Value = -2524765
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H35").Value = 2525000
$ws.Range("I35").Value = 2525000
$ws.Range("K35").Value = 2525000
$ws.Range("M35").Value = -2524710
$ws.Range("H62").Value = 7000
$ws.Range("I62").Value = 7000
$ws.Range("K62").Value = 7000
$ws.Range("M62").Value = -6376
$ws.Range("H65").Value = 7000
$ws.Range("I65").Value = 7000
$ws.Range("K65").Value = 35000
$ws.Range("M65").Value = -31880
$ws.Range("H81").Value = 997.5
$ws.Range("I81").Value = 997.5
$ws.Range("K81").Value = 1995
$ws.Range("M81").Value = -934
$ws.Range("H84").Value = 997.5
$ws.Range("I84").Value = 997.5
$ws.Range("K84").Value = 9975
$ws.Range("M84").Value = -4671
$ws.Range("H132").Value = 1896.1
$ws.Range("I132").Value = 1896.1
$ws.Range("K132").Value = 5688.299999999999
$ws.Range("M132").Value = -3158.299999999999
